$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -124.3497
$ws.Range("B2").Value = -124.3457

$ws.Range("A3").Value = 43.2603
$ws.Range("B3").Value = 43.2632

$ws.Range("A4").Value = -124.2604
$ws.Range("B4").Value = -124.2644

$ws.Range("A5").Value = 43.3219
$ws.Range("B5").Value = 43.319
